$wb = $excel.ActiveWorkbook

# --- Rename Sheet1 -> "ms table" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "ms table"

$ws2 = $wb.Worksheets.Item(2)

# --- Add the new "moved or excluded models" sheet after Sheet2 ---
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "moved or excluded models"

# --- Move rows 7:9 (the DPoRT / Mauritius / Chinese models) from "ms table" to the new sheet ---
$ws1.Range("A7:K9").Copy($ws3.Range("A1:K3"))
$ws3.Range("I1:I3").Clear()

# --- Mark the two still-present "Yes" questionnaire answers on "ms table" ---
$ws1.Range("K5").Value = "Yes"
$ws1.Range("K6").Value = "Yes"

# --- Remove the rows that were moved out ---
$ws1.Rows("7:9").Delete()
